$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Row 52: Select Organization / OrganizationSelect
$row1 = $tbl.ListRows.Add()
$r1 = $row1.Range
$r1.Cells.Item(1,2).Value = "OrganizationSelect"
$r1.Cells.Item(1,1).Value = "Select Organization"
$r1.Cells.Item(1,3).Value = "xpath"
$r1.Cells.Item(1,4).Value = "//select[contains(@title, 'Product Standards')]"
$r1.Cells.Item(1,5).Value = "Step 3"

# Row 53: designation / DesignationInput
$row2 = $tbl.ListRows.Add()
$r2 = $row2.Range
$r2.Cells.Item(1,1).Value = "designation"
$r2.Cells.Item(1,2).Value = "DesignationInput"
$r2.Cells.Item(1,3).Value = "xpath"
$r2.Cells.Item(1,4).Value = "//input[contains(@id, 'Designation')]"
$r2.Cells.Item(1,5).Value = "Step 4"

[void]$ws.Range("A6:D6").Select()
